$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C for rows 2-9 are all zeroed out
$ws.Range("B2:C9").Value = 0

# Column D gets updated specific values
$ws.Range("D2").Value = 0.6505047461123566
$ws.Range("D3").Value = -0.7806995076831166
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = -0.6824410903491035
$ws.Range("D7").Value = 0.7832279335015655
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
